# Applies the data update described in the commit:
# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# Functional change: the "Valor Mora" figures for period 2501 (row 16)
# and period 2412 (row 21) were swapped back to their corrected values.
#   F16: 52000 -> 45066
#   F21: 45066 -> 52000

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 45066
$ws.Range("F21").Value = 52000
